$wb = $excel.ActiveWorkbook

# --- NewLoanInput sheet ---
$wsLoan = $wb.Worksheets.Item("NewLoanInput")

# product value changed (reuses the existing shared-string slot in place)
$wsLoan.Range("B2").Value = "chaithanyatest"

# insert a new "Firstrepaymenton" row after "disbursementon" (adds a brand new shared string)
$wsLoan.Rows("7:7").Insert()
$wsLoan.Range("A7").Value = "Firstrepaymenton"
$wsLoan.Range("B7").Value = 42036

# --- Transactions sheet: update Entry ID value ---
$wsTx = $wb.Worksheets.Item("Transactions")
$wsTx.Range("A2").Value = 207

# --- Update selections on each sheet ---
$wsLoan.Range("B4").Select()

$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B4").Select()

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Range("H8").Select()

# Transactions becomes the active sheet/tab, with C2 selected
$wsTx.Activate()
$wsTx.Range("C2").Select()
